$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 13011.25
$ws.Range("I44").Value = 5045
$ws.Range("K44").Value = 5045
$ws.Range("M44").Value = -4583
$ws.Range("H132").Value = 3057.5454
$ws.Range("I132").Value = 2476.913
$ws.Range("J132").Value = 4393
$ws.Range("K132").Value = 7430.739
$ws.Range("L132").Value = 13179
$ws.Range("M132").Value = -4900.739
$ws.Range("N132").Value = -18239
$ws.Range("H136").Value = 59900
$ws.Range("J136").Value = 59900
$ws.Range("L136").Value = 59900
$ws.Range("N136").Value = -70100
$ws.Range("H137").Value = 1860.129
$ws.Range("I137").Value = 1703.3158
$ws.Range("J137").Value = 2108.4167
$ws.Range("K137").Value = 5109.9474
$ws.Range("L137").Value = 6325.250100000001
$ws.Range("M137").Value = -2559.9474
$ws.Range("N137").Value = -11425.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 84918.086
$ws.Range("I2").Value = 1518.2
$ws.Range("J2").Value = 144489.42
$ws.Range("K2").Value = 1518.2
$ws.Range("L2").Value = 144489.42
$ws.Range("M2").Value = -1405.2
$ws.Range("N2").Value = -144715.42
$ws.Range("H61").Value = 2086.7273
$ws.Range("I61").Value = 2416
$ws.Range("J61").Value = 1510.5
$ws.Range("K61").Value = 2416
$ws.Range("L61").Value = 1510.5
$ws.Range("M61").Value = -2204
$ws.Range("N61").Value = -1934.5
$ws.Range("H74").Value = 891.7646999999999
$ws.Range("I74").Value = 877.3333
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 877.3333
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -3.333300000000008
$ws.Range("N74").Value = -2748
$ws.Range("H77").Value = 891.7646999999999
$ws.Range("I77").Value = 877.3333
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 4386.6665
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -18.66650000000027
$ws.Range("N77").Value = -13736
$ws.Range("H116").Value = 84918.086
$ws.Range("I116").Value = 1518.2
$ws.Range("J116").Value = 144489.42
$ws.Range("K116").Value = 1518.2
$ws.Range("L116").Value = 144489.42
$ws.Range("M116").Value = 775.8
$ws.Range("N116").Value = -149077.42
$ws.Range("H132").Value = 2639.1304
$ws.Range("I132").Value = 2170
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6510
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3980
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 2086.7273
$ws.Range("I136").Value = 2416
$ws.Range("J136").Value = 1510.5
$ws.Range("K136").Value = 7248
$ws.Range("L136").Value = 4531.5
$ws.Range("M136").Value = -4698
$ws.Range("N136").Value = -9631.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 84918.086
$ws.Range("I3").Value = 1518.2
$ws.Range("J3").Value = 144489.42
$ws.Range("K3").Value = 1518.2
$ws.Range("L3").Value = 144489.42
$ws.Range("M3").Value = -1404.2
$ws.Range("N3").Value = -144717.42
$ws.Range("H134").Value = 2401.4878
$ws.Range("I134").Value = 1942.1428
$ws.Range("J134").Value = 3390.8462
$ws.Range("K134").Value = 5826.428400000001
$ws.Range("L134").Value = 10172.5386
$ws.Range("M134").Value = -3291.428400000001
$ws.Range("N134").Value = -15242.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1617.3654
$ws.Range("I31").Value = 1165.5227
$ws.Range("J31").Value = 4102.5
$ws.Range("K31").Value = 1165.5227
$ws.Range("L31").Value = 4102.5
$ws.Range("M31").Value = -870.5227
$ws.Range("N31").Value = -4692.5
$ws.Range("H34").Value = 1617.3654
$ws.Range("I34").Value = 1165.5227
$ws.Range("J34").Value = 4102.5
$ws.Range("K34").Value = 1165.5227
$ws.Range("L34").Value = 4102.5
$ws.Range("M34").Value = -963.5227
$ws.Range("N34").Value = -4506.5
$ws.Range("H58").Value = 883144.0600000001
$ws.Range("I58").Value = 1323668.8
$ws.Range("K58").Value = 1323668.8
$ws.Range("M58").Value = -1323465.8
$ws.Range("H134").Value = 1616.1351
$ws.Range("I134").Value = 1292.6666
$ws.Range("J134").Value = 2213.3076
$ws.Range("K134").Value = 3877.9998
$ws.Range("L134").Value = 6639.9228
$ws.Range("M134").Value = -1342.9998
$ws.Range("N134").Value = -11709.9228
$ws.Range("H136").Value = 883144.0600000001
$ws.Range("I136").Value = 1323668.8
$ws.Range("K136").Value = 3971006.4
$ws.Range("M136").Value = -3968456.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1102.4375
$ws.Range("I5").Value = 1074.3462
$ws.Range("J5").Value = 1224.1666
$ws.Range("K5").Value = 3223.0386
$ws.Range("L5").Value = 3672.4998
$ws.Range("M5").Value = -3111.0386
$ws.Range("N5").Value = -3896.4998
$ws.Range("H92").Value = 519.6667
$ws.Range("I92").Value = 541.6667
$ws.Range("J92").Value = 497.66666
$ws.Range("K92").Value = 1625.0001
$ws.Range("L92").Value = 1492.99998
$ws.Range("M92").Value = -377.0001
$ws.Range("N92").Value = -3988.99998
$ws.Range("H117").Value = 697.7
$ws.Range("J117").Value = 902.25
$ws.Range("L117").Value = 2706.75
$ws.Range("N117").Value = -9590.75
$ws.Range("H131").Value = 11767296
$ws.Range("I131").Value = 9429.166999999999
$ws.Range("J131").Value = 13700096
$ws.Range("K131").Value = 28287.501
$ws.Range("L131").Value = 41100288
$ws.Range("M131").Value = -23247.501
$ws.Range("N131").Value = -41110368
$ws.Range("H132").Value = 2896.6667
$ws.Range("J132").Value = 3495
$ws.Range("L132").Value = 31455
$ws.Range("N132").Value = -36515
$ws.Range("H135").Value = 1102.4375
$ws.Range("I135").Value = 1074.3462
$ws.Range("J135").Value = 1224.1666
$ws.Range("K135").Value = 9669.1158
$ws.Range("L135").Value = 11017.4994
$ws.Range("M135").Value = -7134.1158
$ws.Range("N135").Value = -16087.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2164.35
$ws.Range("I132").Value = 1552.8
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 4658.4
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -2128.4
$ws.Range("N132").Value = -17057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 26900
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 26900
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 26900
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -27240
$ws.Range("H55").Value = 100
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 73
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 69001.336
$ws.Range("I61").Value = 69001.336
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 69001.336
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -68799.336
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 69001.336
$ws.Range("I113").Value = 69001.336
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 69001.336
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -66831.336
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 30000
$ws.Range("J25").Value = 30000
$ws.Range("L25").Value = 30000
$ws.Range("N25").Value = -30586
$ws.Range("H62").Value = 4967.6665
$ws.Range("I62").Value = 4900
$ws.Range("J62").Value = 5001.5
$ws.Range("K62").Value = 4900
$ws.Range("L62").Value = 5001.5
$ws.Range("M62").Value = -4276
$ws.Range("N62").Value = -6249.5
$ws.Range("H65").Value = 4967.6665
$ws.Range("I65").Value = 4900
$ws.Range("J65").Value = 5001.5
$ws.Range("K65").Value = 24500
$ws.Range("L65").Value = 25007.5
$ws.Range("M65").Value = -21380
$ws.Range("N65").Value = -31247.5
$ws.Range("H132").Value = 721.4493
$ws.Range("I132").Value = 586.8644
$ws.Range("J132").Value = 1515.5
$ws.Range("K132").Value = 1760.5932
$ws.Range("L132").Value = 4546.5
$ws.Range("M132").Value = 769.4067999999997
$ws.Range("N132").Value = -9606.5
$ws.Range("H136").Value = 982.1
$ws.Range("I136").Value = 971.0909
$ws.Range("J136").Value = 1012.375
$ws.Range("K136").Value = 2913.2727
$ws.Range("L136").Value = 3037.125
$ws.Range("M136").Value = -363.2727
$ws.Range("N136").Value = -8137.125

Write-Output "Applied 217 cell updates across 8 sheets"